$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63, pushing existing rows 63-75 down to 64-76.
$ws.Rows.Item(63).Insert()

# Populate the new row 63 with the new record.
$ws.Cells.Item(63, 1).Value = 7
$ws.Cells.Item(63, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(63, 3).Value = "Ñuble"
$ws.Cells.Item(63, 4).Value = 44798
$ws.Cells.Item(63, 5).Value = 16
$ws.Cells.Item(63, 6).Value = 100112022
$ws.Cells.Item(63, 7).Value = "Arveja Verde"
$ws.Cells.Item(63, 8).Value = "Perfection"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 60
$ws.Cells.Item(63, 11).Value = 35000
$ws.Cells.Item(63, 12).Value = 36000
$ws.Cells.Item(63, 13).Value = 35500
$ws.Cells.Item(63, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(63, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(63, 16).Value = 1420
$ws.Cells.Item(63, 17).Value = 25
$ws.Cells.Item(63, 18).Value = "Hortaliza"
